$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 0. Materialize column I (new, blank, trailing column) so the used range
#    extends to I1:I6 like the target. Setting NumberFormat on an empty cell
#    is enough to make the engine persist it without touching its value.
# ---------------------------------------------------------------------------
$ws.Range("I1:I6").NumberFormat = "General"

# ---------------------------------------------------------------------------
# 1. Clear the old "Residuals" row content first -- while the row still
#    carries its original style, so the (now-blank) row survives as an
#    empty-but-present row instead of being dropped from sheetData.
# ---------------------------------------------------------------------------
$ws.Range("A6:I6").ClearContents()

# ---------------------------------------------------------------------------
# 2. Re-label row 2 headers: the ANOVA "F" / "P" columns became the
#    Anova()-style "Chisq" / "Pr(>Chisq)" columns.
# ---------------------------------------------------------------------------
$ws.Range("C2").Value = "Chisq"
$ws.Range("D2").Value = "Pr(>Chisq)"
$ws.Range("E2").Value = "Chisq"
$ws.Range("F2").Value = "Pr(>Chisq)"
$ws.Range("G2").Value = "Chisq"
$ws.Range("H2").Value = "Pr(>Chisq)"

# ---------------------------------------------------------------------------
# 3. Refresh the statistics for each model term with the re-run numbers.
# ---------------------------------------------------------------------------
# Elevation
$ws.Range("C3").Value = 0.14799092769068201
$ws.Range("D3").Value = 0.70046274069761305
$ws.Range("E3").Value = 1.10383322144877005
$ws.Range("F3").Value = 0.29342640883534998
$ws.Range("G3").Value = 1.63923606406755007
$ws.Range("H3").Value = 0.20043031540351400

# Fire
$ws.Range("C4").Value = 4.05261052790720999
$ws.Range("D4").Value = 0.04410309356358850
$ws.Range("E4").Value = 0.12366749547507599
$ws.Range("F4").Value = 0.72509034022772600
$ws.Range("G4").Value = 0.42510317744731801
$ws.Range("H4").Value = 0.51440157932468500

# Elevation*Fire
$ws.Range("C5").Value = 0.000548621287667824
$ws.Range("D5").Value = 0.98131312502373802
$ws.Range("E5").Value = 0.39190331494243502
$ws.Range("F5").Value = 0.53130063162199404
$ws.Range("G5").Value = 1.70687135934074008
$ws.Range("H5").Value = 0.19139179930618100

# ---------------------------------------------------------------------------
# 4. Drop the small 10pt font + thin box border that used to decorate the
#    table -- the refreshed table uses plain default formatting instead.
# ---------------------------------------------------------------------------
$ws.Range("A1:I6").ClearFormats()

# ---------------------------------------------------------------------------
# 5. Every numeric statistic (Df / Chisq / Pr) now shares the "0.000"
#    number format, instead of only half the columns.
# ---------------------------------------------------------------------------
$ws.Range("B3:I5").NumberFormat = "0.000"

# ---------------------------------------------------------------------------
# 6. Resize the columns to fit the new header/value text.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 11.998697916666666
$ws.Columns.Item(2).ColumnWidth = 4.830729166666667
$ws.Columns.Item(3).ColumnWidth = 6.498697916666667
$ws.Columns.Item(4).ColumnWidth = 8.666666666666666
$ws.Columns.Item(5).ColumnWidth = 6.830729166666667
$ws.Columns.Item(6).ColumnWidth = 8.666666666666666
$ws.Columns.Item(7).ColumnWidth = 7.830729166666667
$ws.Columns.Item(8).ColumnWidth = 8.666666666666666

# ---------------------------------------------------------------------------
# 7. Match the author's on-screen view: zoomed to 150%, selection on the
#    block of refreshed statistics.
# ---------------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.Zoom = 150
$ws.Range("C3:H5").Select()
